# Apply the "Hjemme passive tweaks lichtwark deleted values" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (subject id headers) - B1:E1 changed
$ws.Range("B1").Value2 = 16
$ws.Range("C1").Value2 = 20
$ws.Range("D1").Value2 = 16
$ws.Range("E1").Value2 = 20

# Row 2 (CON) - B2, D2, E2 values removed entirely; C2 replaced with new value
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value2 = -2.9033381510991703
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 (STR) - B3:E3 replaced with new values
$ws.Range("B3").Value2 = -3.105531684919832
$ws.Range("C3").Value2 = 3.9959297561476745
$ws.Range("D3").Value2 = -0.39689215022412583
$ws.Range("E3").Value2 = 11.749425093518212

# Match the updated selection in the sheet view
$ws.Range("B1:E3").Select()
